# Fruta / hortaliza, semanal
# Update weekly price records for "Achicoria" in the Lo Valledor wholesale
# market sheet: shift/update existing rows 5-14 and append two new rows
# (15 and 16) with additional weekly observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($Row, $Values)
    foreach ($col in $Values.Keys) {
        $colIndex = [int][char]$col - [int][char]'A' + 1
        $ws.Cells.Item($Row, $colIndex).Value = $Values[$col]
    }
}

# --- Row 5: update date, volume, prices and origin ---
Set-RowValues 5 @{
    'D' = 44236
    'J' = 180
    'K' = 4000
    'L' = 4500
    'M' = 4167
    'O' = 'Región Metropolitana'
    'P' = 260
}

# --- Row 6: update date and volume ---
Set-RowValues 6 @{
    'D' = 44186
    'J' = 160
}

# --- Row 7: update date and volume ---
Set-RowValues 7 @{
    'D' = 44188
    'J' = 210
}

# --- Row 8: update date ---
Set-RowValues 8 @{
    'D' = 44232
}

# --- Row 9: update date, max/avg price and $/Kg ---
Set-RowValues 9 @{
    'D' = 44846
    'L' = 5000
    'M' = 5000
    'P' = 312
}

# --- Row 10: update date and volume ---
Set-RowValues 10 @{
    'D' = 44189
    'J' = 250
}

# --- Row 11: update date ---
Set-RowValues 11 @{
    'D' = 44215
}

# --- Row 12: update date and volume ---
Set-RowValues 12 @{
    'D' = 44210
    'J' = 340
}

# --- Row 13: update date and volume ---
Set-RowValues 13 @{
    'D' = 44231
    'J' = 250
}

# --- Row 14: update date, volume, min price, avg price, origin and $/Kg ---
Set-RowValues 14 @{
    'D' = 44204
    'J' = 430
    'K' = 5000
    'M' = 5500
    'O' = 'Provincia de Quillota'
    'P' = 344
}

# --- Row 15: new record ---
Set-RowValues 15 @{
    'A' = 6
    'B' = 'Mercado Mayorista Lo Valledor de Santiago'
    'C' = 'Metropolitana'
    'D' = 44187
    'E' = 13
    'F' = 100112010
    'G' = 'Achicoria'
    'H' = 'Sin especificar'
    'I' = 'Primera'
    'J' = 160
    'K' = 5000
    'L' = 6000
    'M' = 5500
    'N' = '$/caja 16 unidades'
    'O' = 'Provincia de Quillota'
    'P' = 344
    'Q' = 16
    'R' = 'Hortaliza'
}

# --- Row 16: new record ---
Set-RowValues 16 @{
    'A' = 6
    'B' = 'Mercado Mayorista Lo Valledor de Santiago'
    'C' = 'Metropolitana'
    'D' = 44292
    'E' = 13
    'F' = 100112010
    'G' = 'Achicoria'
    'H' = 'Sin especificar'
    'I' = 'Primera'
    'J' = 90
    'K' = 6000
    'L' = 6000
    'M' = 6000
    'N' = '$/caja 16 unidades'
    'O' = 'Región Metropolitana'
    'P' = 375
    'Q' = 16
    'R' = 'Hortaliza'
}

# Ensure the date columns on the two new rows use the same number format
# (date/time) as the rest of column D.
$ws.Cells.Item(15, 4).NumberFormat = $ws.Cells.Item(14, 4).NumberFormat
$ws.Cells.Item(16, 4).NumberFormat = $ws.Cells.Item(14, 4).NumberFormat

Write-Output "Applied Achicoria weekly update: rows 5-14 revised, rows 15-16 added."
